$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. "20.13")
# must be forced to stay as text, matching the source inlineStr cells,
# otherwise Excel auto-converts them into numeric values.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "26.991.17"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "1.658.50"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "215.31"
$ws.Range("E5").Value = "  +1.43%  "
Set-TextValue $ws.Range("D6") "0.508"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("E7").Value = "  -0.10%  "
Set-TextValue $ws.Range("D8") "0.251"
$ws.Range("E8").Value = "  +2.62%  "
$ws.Range("E9").Value = "  +2.03%  "
Set-TextValue $ws.Range("D10") "20.13"
$ws.Range("E10").Value = "  +4.76%  "
Set-TextValue $ws.Range("D11") "0.0883"
$ws.Range("E11").Value = "  +4.28%  "
$ws.Range("D12").Value = "1.892.07"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D13").Value = "1.656.92"
$ws.Range("E13").Value = "  +2.77%  "
Set-TextValue $ws.Range("D14") "4.09"
$ws.Range("E14").Value = "  +2.12%  "
Set-TextValue $ws.Range("D15") "0.523"
$ws.Range("E15").Value = "  +2.71%  "
Set-TextValue $ws.Range("D16") "65.57"
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "26.987.19"
$ws.Range("E17").Value = "  +2.19%  "
Set-TextValue $ws.Range("D18") "236.56"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "0.0₃0738"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +3.93%  "
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("E24").Value = "  +2.62%  "
Set-TextValue $ws.Range("D25") "145.20"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("E27").Value = "  +0.74%  "
Set-TextValue $ws.Range("D28") "15.85"
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("D32").Value = "1.560.69"
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("E34").Value = "  +4.74%  "
$ws.Range("E35").Value = "  +8.36%  "
Set-TextValue $ws.Range("D36") "2.41"
$ws.Range("E36").Value = "  -0.39%  "
Set-TextValue $ws.Range("D37") "0.580"
$ws.Range("E37").Value = "  +3.30%  "
Set-TextValue $ws.Range("D38") "0.902"
$ws.Range("E38").Value = "  +9.15%  "
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("E41").Value = "  -0.07%  "
Set-TextValue $ws.Range("D42") "66.48"
$ws.Range("E42").Value = "  +8.44%  "
$ws.Range("E43").Value = "  +6.27%  "
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("D45").Value = "1.800.85"
$ws.Range("E45").Value = "  +2.91%  "
Set-TextValue $ws.Range("D46") "0.774"
$ws.Range("E46").Value = "  +1.65%  "
Set-TextValue $ws.Range("D47") "90.20"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("E48").Value = "  +2.96%  "
Set-TextValue $ws.Range("D49") "0.1000"
$ws.Range("E49").Value = "  +4.37%  "
$ws.Range("E50").Value = "  +0.98%  "
Set-TextValue $ws.Range("D51") "7.70"
$ws.Range("E51").Value = "  +2.96%  "
